$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A2 to the MATCH query text (new shared string). Single-quoted PowerShell
# string so the embedded backticks are treated literally (not as escape chars).
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Bone cancer, NOS''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# Row 2 grows to fit the wrapped query text (ht="87" in the target sheet).
$ws.Rows.Item(2).RowHeight = 87

# Selection moves from C7 to the A2:A6 block (used by the next step of the
# automation to highlight/copy the pasted query).
$ws.Range("A2:A6").Select()
